$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Metadata")

# Translate Publisher value (row 9) from German to English
$ws.Range("B9").Value = "Independent Trusted Third Party of the University Medicine Greifswald"

# Translate Contact value (row 10) from German to English
$ws.Range("B10").Value = "Independent Trusted Third Party of the University Medicine Greifswald (https://www.ths-greifswald.de/)"

# Add a Description value (row 12), previously empty
$ws.Range("B12").Value = "Element names for use in `$requestTasks operation. "
